# Append a new annotation row (row 6) for "parisk", mirroring the existing
# rows 2-5 in the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "parisk"
$ws.Range("B6").Value = 3

# polite_expressions (column C) is blank for every row in this sheet, stored
# as an empty *text* cell (not a numeric blank). A leading apostrophe forces
# text-empty-string semantics; then copy the plain (non quote-prefixed)
# style from the row above so no stray "treat as text" formatting sticks.
$ws.Range("C6").Value = "'"
$ws.Range("C6").Style = $ws.Range("C5").Style

$ws.Range("D6").Value = "ACK"
$ws.Range("E6").Value = "EXP"
$ws.Range("F6").Value = "42be9703-0e9b-4ce8-962d-60bf1f233ce8"
$ws.Range("G6").Value = "SJCPLLpaW_annotated.xlsx"
$ws.Range("H6").Value = "The results show that DeePa achieves speedups compared to PyTorch and TensorFlow with all of the tested minibatch sizes."
